$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(13.46456841828631, 7.237647162246065, 14.00369035395992, 47.40230617507002, 3.711498604867593, 11.29078040900474, 9.595067975462907, 18.99769123883691, 36.86143484149508)
    3 = @(13.44089697468744, 7.249130751239073, 14.02303258333937, 47.12862184932456, 3.715841198117702, 11.31216645505469, 9.615335707911068, 18.40031508502705, 36.66990447118214)
    4 = @(13.42938265462406, 7.256736956642992, 14.03716565015127, 46.97262192089095, 3.718646908430661, 11.32691093684996, 9.628643202279045, 18.02485520896364, 36.56157404316777)
    5 = @(13.42545186324604, 7.259976432875316, 14.04349205974699, 46.91212518869624, 3.719825436473566, 11.33332499675198, 9.634283630441059, 17.86990355188765, 36.51978946809245)
    6 = @(13.42484518937866, 7.260522802078719, 14.04457679044532, 46.90226669068161, 3.720023258950126, 11.33441453733451, 9.635233371458524, 17.84406337566597, 36.51299456900767)
    7 = @(13.42932655763089, 7.256780078595019, 14.03724867492417, 46.97179353550058, 3.718662659852912, 11.32699579719412, 9.628718389690167, 18.02277304767602, 36.56100092582795)
    8 = @(13.45578106520091, 7.241491624922625, 14.00989090444787, 47.30546942562952, 3.712967095336765, 11.2978193922245, 9.601877470073127, 18.79364780656866, 36.79349024983726)
    9 = @(13.53150781846182, 7.215904769877893, 13.97416704116737, 48.05310487567979, 3.702897338804882, 11.25341156141802, 9.556067829806155, 20.2273683202997, 37.32147138025261)
    10 = @(13.60149387191252, 7.199768752099516, 13.95886894994027, 48.65592705651245, 3.696160297802366, 11.22860026253219, 9.526541753109306, 21.22223697909767, 37.75110129544031)
    11 = @(13.63639464457916, 7.193002898745464, 13.95428940830525, 48.9409917901105, 3.693237113820094, 11.21901114985897, 9.513999979181243, 21.66018057919901, 37.95507418711843)
    12 = @(13.65004563244606, 7.190523199770567, 13.95289744902877, 49.05042410277704, 3.692150384744215, 11.21562416779603, 9.509378190022858, 21.82377585682186, 38.03349089720398)
    13 = @(13.64708640699695, 7.191053586948757, 13.95318201306632, 49.02679108287216, 3.692383534264153, 11.21634275281289, 9.510367911024268, 21.78864458690801, 38.01655095156429)
    14 = @(13.63750904012103, 7.192797242649865, 13.9541680327443, 48.9499655354394, 3.693147303550136, 11.21872760616853, 9.51361718886073, 21.67368539489659, 37.96150229263616)
    15 = @(13.63169907231224, 7.193876004016351, 13.95481656350242, 48.90309874599514, 3.693617763188286, 11.22022020260255, 9.515624057686194, 21.60297336126124, 37.92793509483962)
    16 = @(13.59927419815448, 7.20022245727421, 13.95921612271155, 48.63750923685408, 3.696354171180194, 11.22926110294041, 9.527379254438854, 21.19330956972086, 37.73793862966065)
    17 = @(13.58016340057848, 7.204262780214568, 13.96252464355456, 48.4773028658658, 3.698069022147965, 11.23524228606841, 9.534818266726242, 20.93814219015166, 37.62353220979354)
    18 = @(13.56946013472474, 7.206640755558611, 13.96465160525664, 48.3861811525388, 3.699068688541915, 11.23884229425478, 9.53918076689448, 20.79000725568362, 37.55853491760423)
    19 = @(13.56588595247156, 7.20745519393721, 13.96541022556537, 48.35550711223603, 3.699409451859707, 11.24008863616215, 9.540672236809419, 20.73962067985785, 37.53666791909159)
    20 = @(13.58216793176924, 7.203827084665053, 13.96214926343796, 48.49425157576233, 3.697885094761882, 11.23458904057303, 9.534017703779316, 20.96544799483449, 37.63562788101501)
    21 = @(13.64031039274446, 7.19228285467749, 13.95386912761644, 48.97249137689764, 3.69292241837717, 11.21802048871065, 9.512659340151639, 21.70751365554064, 37.97763989797382)
    22 = @(13.68084124487569, 7.185218118467729, 13.9504521990589, 49.29366513159238, 3.689796801958936, 11.20861542338142, 9.499443429866417, 22.17935961385674, 38.20799774611488)
    23 = @(13.65897959643782, 7.188944848577463, 13.95209339001096, 49.12148544161741, 3.691454269784962, 11.21350481907321, 9.506429168724896, 21.92877110911181, 38.08444339419763)
    24 = @(13.58126079885887, 7.204023890857137, 13.9623182724373, 48.48658599641973, 3.697968205414837, 11.23488387031035, 9.534379371277478, 20.95310750188673, 37.63015700569962)
    25 = @(13.50848548590518, 7.222357963615265, 13.98190979126697, 47.84118917545405, 3.705504717520317, 11.26405308227006, 9.567733037336, 19.84905939529497, 37.17113856588063)
}

$cols = @("C", "D", "E", "F", "G", "J", "L", "N", "O")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $row
        $ws.Range($addr).Value2 = $vals[$i]
    }
}
